# Update column F (dSF) values for specific rows to reflect
# the repulled data / recalculated mean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    4  = 2
    7  = -5
    9  = -1
    10 = -8
    11 = -12
    13 = 8
    14 = -5
    15 = -4
    18 = -2
    19 = -4
    20 = -1
    23 = -10
    24 = -2
    25 = 7
    27 = 10
    29 = -4
    34 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
